$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for price (D) column cells so numeric-looking strings are preserved exactly
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D11", "D12", "D13", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D33", "D34", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply cell value changes
$ws.Range("D2").Value = '66.643.88'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").Value = '3.551.19'
$ws.Range("E3").Value = '  -1.30%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '608.14'
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("D6").Value = '144.61'
$ws.Range("E6").Value = '  -2.80%  '
$ws.Range("D7").Value = '3.553.61'
$ws.Range("E7").Value = '  -1.11%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  +5.73%  '
$ws.Range("E10").Value = '  -2.77%  '
$ws.Range("D11").Value = '7.79'
$ws.Range("E11").Value = '  -3.43%  '
$ws.Range("D12").Value = '0.413'
$ws.Range("E12").Value = '  -0.81%  '
$ws.Range("D13").Value = '4.160.25'
$ws.Range("E13").Value = '  -1.04%  '
$ws.Range("E14").Value = '  -5.95%  '
$ws.Range("D15").Value = '29.02'
$ws.Range("E15").Value = '  -2.78%  '
$ws.Range("D16").Value = '3.552.69'
$ws.Range("E16").Value = '  -1.48%  '
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").Value = '0.117'
$ws.Range("E17").Value = '  +0.77%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '66.604.21'
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").Value = '10.94'
$ws.Range("E19").Value = '  -4.93%  '
$ws.Range("E20").Value = '  -1.88%  '
$ws.Range("D21").Value = '14.71'
$ws.Range("E21").Value = '  -2.67%  '
$ws.Range("D22").Value = '426.99'
$ws.Range("E22").Value = '  -0.34%  '
$ws.Range("D23").Value = '0.598'
$ws.Range("E23").Value = '  -3.31%  '
$ws.Range("D24").Value = '77.60'
$ws.Range("E24").Value = '  -1.61%  '
$ws.Range("D25").Value = '3.697.23'
$ws.Range("E25").Value = '  -1.23%  '
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("D27").Value = '0.0000116'
$ws.Range("E27").Value = '  -5.51%  '
$ws.Range("D28").Value = '8.02'
$ws.Range("E28").Value = '  -3.79%  '
$ws.Range("D29").Value = '2.48'
$ws.Range("E29").Value = '  -1.69%  '
$ws.Range("D30").Value = '9.04'
$ws.Range("E30").Value = '  -5.45%  '
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D33").Value = '3.562.44'
$ws.Range("E33").Value = '  -0.83%  '
$ws.Range("D34").Value = '24.49'
$ws.Range("E34").Value = '  -3.84%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("E36").Value = '  -8.02%  '
$ws.Range("D37").Value = '7.63'
$ws.Range("E37").Value = '  -3.08%  '
$ws.Range("E38").Value = '  -4.32%  '
$ws.Range("D39").Value = '177.00'
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("D40").Value = '5.30'
$ws.Range("E40").Value = '  -6.32%  '
$ws.Range("D41").Value = '0.0829'
$ws.Range("E41").Value = '  -3.35%  '
$ws.Range("D42").Value = '5.03'
$ws.Range("E42").Value = '  -4.26%  '
$ws.Range("D43").Value = '0.864'
$ws.Range("E43").Value = '  -3.82%  '
$ws.Range("D44").Value = '45.34'
$ws.Range("E44").Value = '  -1.82%  '
$ws.Range("E45").Value = '  -6.55%  '
$ws.Range("E46").Value = '  +0.18%  '
$ws.Range("D47").Value = '2.40'
$ws.Range("E47").Value = '  -6.45%  '
$ws.Range("D48").Value = '7.14'
$ws.Range("E48").Value = '  -0.76%  '
$ws.Range("D49").Value = '23.36'
$ws.Range("E49").Value = '  -2.82%  '
$ws.Range("D50").Value = '1.12'
$ws.Range("E50").Value = '  -4.96%  '
$ws.Range("D51").Value = '0.921'
$ws.Range("E51").Value = '  -3.47%  '

